$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has a pre-existing quirk: column A holds field labels, while
# columns B/C hold the corresponding values offset by a couple of rows.
# Row 13 previously held only B/C data (no A label) for "Docentes
# responsaveis:". That row is removed, shifting everything below up by one,
# and a handful of B/C values are corrected/re-pointed afterwards.

$ws.Rows("13:13").Delete()

# Fix the non-breaking space in the English name value (B4/C4).
$ws.Range("B4:C4").Value = "Case study in environmental impact"

# Row 10 (Objetivos:) now shows the professor's name value.
$ws.Range("B10:C10").Value = "5840938 - Marcelo Rodrigues de Holanda"

# Row 13 (Programa resumido:) now shows "Semestral".
$ws.Range("B13:C13").Value = "Semestral"

# Row 15 (Programa:) now shows the activation date. Copy the existing cell
# (with its formatting) from B8/C8 instead of typing the literal, so Excel
# does not reinterpret the date-shaped string as a numeric date serial.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# Row 18 (Metodo:) now shows the professor's name value.
$ws.Range("B18:C18").Value = "5840938 - Marcelo Rodrigues de Holanda"

# Row 19 (Criterio:) now shows the lecture method text.
$ws.Range("B19:C19").Value = "Aula expositiva e exercícios dirigidos."

# Row 20 (Norma de recuperacao:) now shows the weighted-average text.
$ws.Range("B20:C20").Value = "Média ponderada de exercícios e provas."

# Row 21 (Bibliografia:) now shows the passing-grade text.
$ws.Range("B21:C21").Value = "Prova única com nota igual ou superior a 5,0."
